# Update countries & provincias Spain
# Applies updated COVID-19 stats for a handful of countries (rows identified by
# the "Pais" column in worksheet "Pais"): Estados Unidos (row 4), Turquia (row 10),
# Canada (row 15), and Marruecos (row 55).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 995288
$ws.Range("C4").Value = 8128
$ws.Range("D4").Value = 122675
$ws.Range("E4").Value = 816474
$ws.Range("G4").Value = 726
$ws.Range("H4").Value = 56139

# Row 10 - Turquia
$ws.Range("B10").Value = 112261
$ws.Range("C10").Value = 2131
$ws.Range("D10").Value = 33791
$ws.Range("E10").Value = 75570
$ws.Range("F10").Value = 1736
$ws.Range("G10").Value = 95
$ws.Range("H10").Value = 2900

# Row 15 - Canada
$ws.Range("B15").Value = 47346
$ws.Range("C15").Value = 451
$ws.Range("D15").Value = 17916
$ws.Range("E15").Value = 26813

# Row 55 - Marruecos
$ws.Range("B55").Value = 4120
$ws.Range("C55").Value = 55
$ws.Range("D55").Value = 695
$ws.Range("E55").Value = 3263
$ws.Range("G55").Value = 1
$ws.Range("H55").Value = 162
